$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Italy")
$ws.Copy([System.Reflection.Missing]::Value, $ws)

$newWs = $wb.Worksheets.Item($ws.Index + 1)
$newWs.Name = "Spain"
$newWs.Range("B2").Value = "Spain Market"
$newWs.Range("B4").Value = "NGC-3103/T2064/T2063"

$newWs.Columns.AutoFit() | Out-Null

$ws.Range("A1:D12").Select()

$newWs.Activate()
$newWs.Range("E9").Select()
